$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active/selected tab (was NewLoanInput before)
$ws.Activate()

# Insert a new blank column before the old column N (Late/Outstanding block),
# shifting N->O, O->P, P->Q
$ws.Columns("N").Insert()

# Approximate the original column width for the newly inserted column
$ws.Columns("N").ColumnWidth = 9.83

# Update the selected cell on the Repayment schedule sheet
$ws.Range("S7").Select() | Out-Null
